# Apply "Updated Matched, Files from s3" reconciliation refresh.
#
# Summary of the change:
#  - matched: the two duplicated ACH payment rows are removed; sheet reverts
#    to the empty "info"/"No data" placeholder state.
#  - unmatched_invoices: gains the two outstanding invoices that no longer
#    match anything (date/description/amount).
#  - unmatched_payments: gains the two ACH payments (now considered
#    unmatched bank-statement lines) inserted at the top, with the running
#    "balance" column recalculated for every subsequent row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: matched  ->  reset back to "No data" placeholder
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("matched")

$ws1.Range("B1:C5").Clear()
$ws1.Range("A3:A5").Clear()

$ws1.Range("A1").Value = "info"
$ws1.Range("A2").Value = "No data"

# ---------------------------------------------------------------------
# Sheet: unmatched_invoices  ->  populate with the two unmatched invoices
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("unmatched_invoices")

# Extend the header style (currently only on A1) across B1:C1.
$ws2.Range("A1").Copy($ws2.Range("B1"))
$ws2.Range("A1").Copy($ws2.Range("C1"))
$ws2.Range("A1").Value = "date"
$ws2.Range("B1").Value = "description"
$ws2.Range("C1").Value = "amount"

$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "2025-06-23"
$ws2.Range("A2").ClearFormats()
$ws2.Range("B2").Value = "Invoice INV-20250623-59B50E89"
$ws2.Range("C2").Value = 28518.85

$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "2025-06-23"
$ws2.Range("A3").ClearFormats()
$ws2.Range("B3").Value = "Invoice INV-20250623-9406A583"
$ws2.Range("C3").Value = 27033.29

# ---------------------------------------------------------------------
# Sheet: unmatched_payments  ->  insert the two ACH payments at the top
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("unmatched_payments")

$ws3.Range("A2:E3").Insert()
$ws3.Range("A2:E3").ClearFormats()

$ws3.Range("A2").NumberFormat = "@"
$ws3.Range("A2").Value = "2025-06-21"
$ws3.Range("A2").ClearFormats()
$ws3.Range("B2").Value = "ACH Payment - Global Mfg Corp INV-20250620-D0E0F066"
$ws3.Range("C2").Value = 28518.85
$ws3.Range("D2").Value = 153518.85
$ws3.Range("E2").Value = "credit"

$ws3.Range("A3").NumberFormat = "@"
$ws3.Range("A3").Value = "2025-06-22"
$ws3.Range("A3").ClearFormats()
$ws3.Range("B3").Value = "ACH Payment - Omkar Mestry INV-20250620-996A7766"
$ws3.Range("C3").Value = 27033.29
$ws3.Range("D3").Value = 180552.14
$ws3.Range("E3").Value = "credit"

# Recalculate the running balance for the rows that followed (now rows 4-10)
$ws3.Range("D4").Value = 176552.14
$ws3.Range("D5").Value = 176362.15
$ws3.Range("D6").Value = 176282.16
$ws3.Range("D7").Value = 121282.16
$ws3.Range("D8").Value = 120532.16
$ws3.Range("D9").Value = 139282.16
$ws3.Range("D10").Value = 138961.66

Write-Host "Reconciliation sheets refreshed from S3 match results."
